$d = $word.ActiveDocument

$d.Content.Find.Execute("Class Diagram", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Activity Diagram", 2)
